# Weights.xlsx: refresh the estimated-mass/weight results after an update
# to the underlying Aircraft-class weight estimation (JPADCore_v2).
# Each worksheet holds a small table of method-by-method estimates in
# column C (and, where present, the percent-error vs. reference in
# column D); this script rewrites the recomputed values in place.

$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS: updated Mass rows (Maximum TO/TO/Landing Mass,
#     Maximum Zero Fuel/Zero Fuel/Operating Empty/Empty Mass) and the
#     mirrored Weight rows in Newtons.
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C6").Value = 50991.47705740609
$ws.Range("C7").Value = 50991.47705740609
$ws.Range("C8").Value = 45892.32935166548
$ws.Range("C12").Value = 39032.41037812544
$ws.Range("C13").Value = 39032.41037812544
$ws.Range("C14").Value = 26162.41037812544
$ws.Range("C15").Value = 25433.323087125438
$ws.Range("C16").Value = 24583.113087125446
$ws.Range("C20").Value = 500055.56848501135
$ws.Range("C21").Value = 500055.56848501135
$ws.Range("C22").Value = 450050.01163651014
$ws.Range("C26").Value = 382777.18723464373
$ws.Range("C27").Value = 382777.18723464373
$ws.Range("C28").Value = 256565.6017346438
$ws.Range("C29").Value = 249415.69785235863
$ws.Range("C30").Value = 241077.98595585872

# --- FUSELAGE: NICOLAI_1984, RAYMER, ROSKAM methods + the averaged
#     "Estimated Mass" row.
$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C7").Value = 5694.0
$ws.Range("D7").Value = 13.77305332986995
$ws.Range("C8").Value = 5645.0
$ws.Range("D8").Value = 12.793973664755155
$ws.Range("C9").Value = 5934.0
$ws.Range("D9").Value = 18.56854556716689
$ws.Range("C12").Value = 6386.833333333332
$ws.Range("D12").Value = 27.61670696212227

# --- WING: all comparison methods (JENKINSON .. KROO) + "Estimated Mass".
$ws = $wb.Worksheets.Item("WING")
$ws.Range("C7").Value = 4580.0
$ws.Range("D7").Value = 42.71914243869
$ws.Range("C8").Value = 3645.0
$ws.Range("D8").Value = 13.583247639525121
$ws.Range("C9").Value = 4676.0
$ws.Range("D9").Value = 45.71063538063634
$ws.Range("C10").Value = 4407.0
$ws.Range("D10").Value = 37.3282228662242
$ws.Range("C11").Value = 4775.0
$ws.Range("D11").Value = 48.79561247701851
$ws.Range("C12").Value = 4214.0
$ws.Range("D12").Value = 31.314075597519576
$ws.Range("C13").Value = 3756.7142857142853
$ws.Range("D13").Value = 17.064419485659087

# --- HORIZONTAL TAIL: RAYMER method + "Estimated Mass".
$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C9").Value = 437.0
$ws.Range("D9").Value = -20.299106328652194
$ws.Range("C10").Value = 467.0
$ws.Range("D10").Value = -14.827649097209544

# --- VERTICAL TAIL: SADRAY method + "Estimated Mass".
$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C8").Value = 412.0
$ws.Range("D8").Value = -8.990501435829454
$ws.Range("C9").Value = 475.5
$ws.Range("D9").Value = 5.036447978793905

# --- LANDING GEARS: all comparison methods (JENKINSON .. TORENBEEK_2013)
#     + "Estimated Mass".
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C5").Value = 1488.0
$ws.Range("D5").Value = -33.36617258519545
$ws.Range("C6").Value = 2036.0
$ws.Range("D6").Value = -8.82629528458195
$ws.Range("C7").Value = 2315.0
$ws.Range("D7").Value = 3.667547355693904
$ws.Range("C8").Value = 2006.0
$ws.Range("D8").Value = -10.169719224396557
$ws.Range("C9").Value = 1961.25
$ws.Range("D9").Value = -12.173659934620025
